$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update threshold values ---
# alpha_distance_range: Min 5.6 -> 5
$ws.Range("B2").Value = 5
# beta_distance_range: Min 5.6 -> 4.5
$ws.Range("B3").Value = 4.5
# ratio_threshold_range: Max 1.3 -> 1.5
$ws.Range("C4").Value = 1.5

# --- Remove the theta_threshold_range row entirely (was row 5) ---
# this shifts the old row 6 (pie_threshold_range) up to row 5
$ws.Rows(5).Delete()

# --- Update the (now-shifted) pie_threshold_range row ---
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 20

# --- Narrow column C's best-fit width now that the longest value (172.3) is gone ---
$ws.Columns("C").ColumnWidth = 4.8125

# --- Match the new selection left behind on the sheet ---
$ws.Range("C5").Select()

# --- Page setup now present on the sheet ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
